$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebalanced lvl 04 -----------------------------------------------------
# Tile "03-Curva-Baixo-Dir-01" (row 11) now also appears once in "lvl 6"
# (column I).
$ws.Range("I11").Value = 1

# Re-apply the row-22 "Tempo minimo gasto por lvl" formula across the whole
# D22:J22 range in one shot so Excel collapses it back into a single shared
# formula group (matches the workbook's usual style) and recalculates the
# dependent totals (I21, I22, O23) with the updated lvl 6 count.
$ws.Range("D22:J22").Formula = "=(D4*`$M`$4+D5*`$M`$5+D6*`$M`$6+D7*`$M`$7+D8*`$M`$8+D9*`$M`$9+D10*`$M`$10+D11*`$M`$11+D12*`$M`$12+D13*`$M`$13+D14*`$M`$14+D15*`$M`$15+D16*`$M`$16+D17*`$M`$17+D18*`$M`$18+D19*`$M`$19+5)/D21*D20"

# --- Fix Victory Menu -------------------------------------------------------
# Scroll the sheet so row 9 is at the top and select I13, matching the
# reopened view of the fixed Victory Menu layout.
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I13").Select()

$wb.Save()
